$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date stamp that is identical for every
# data row (rows 2-115). This automatic update advances that date by one
# day (serial date 46082 -> 46083) for every row.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $current = $cell.Value2
    if ($current -ne $null) {
        $cell.Value = $current + 1
    }
}
